$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting Features/Adjust Factor/etc. right.
$ws.Columns.Item(1).Insert()

# New header for the inserted "Priority" column (bold, matching the other headers in row 2).
$ws.Range("A2").Value = "Priority"
$ws.Range("A2").Font.Bold = $true

# Sequential priority numbers for each backlog row (rows 3-18).
for ($r = 3; $r -le 18; $r++) {
    $ws.Cells.Item($r, 1).Value = ($r - 2)
}

# Selection moves to A2 after the edit.
$ws.Range("A2").Select()
